$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'308.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.13%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'40.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'2.06%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.114"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.14%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07626"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.58%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.607"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.13%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.9031"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.33%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D9").Value = "'0.1112"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'9.35%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1782"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.11%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09097"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.08%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04204"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-5.22%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1051"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.47%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001257"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.00%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005746"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.97%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.349"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.18%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'4.243"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.50%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'0.50%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.574"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-7.26%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'1.79%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.2830"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.51%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04064"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-2.84%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001229"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.96%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004106"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.76%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'-0.11%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003748"
$ws.Range("D26").Style = "Normal"
$ws.Range("D38").Value = "'0.02408"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'1.37%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05174"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-1.01%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007758"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-2.36%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1304"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.73%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007050"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'11.22%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.001950"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-0.46%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008796"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'0.30%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3078"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-7.86%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006953"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'5.71%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.02%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.03085"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'696.50%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.02%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.02%"
$ws.Range("E51").Style = "Normal"